# "100 smartphone, updated rankings"
# Extend the rankings table on Sheet1 from 25 data rows (A2:C26) to
# 100 data rows (A2:C101):
#   - rows 27-63: score 106 / rank -1482 (same as row 26) with the date
#     series continuing day-by-day
#   - row 64: score bumps to 129 / rank -1326 (the "100th smartphone"
#     milestone), date 44368
#   - rows 65-101: only the date column is known so far (score/rank not
#     filled in yet), dates continuing day-by-day through 44405
# Also keep the line chart's series ranges in sync with the new extent.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 27-63: continue the existing 106 / -1482 plateau, one day at a time.
for ($r = 27; $r -le 63; $r++) {
    $date = 44330 + ($r - 26)
    $ws.Cells.Item($r, 1).Value = $date
    $ws.Cells.Item($r, 2).Value = 106
    $ws.Cells.Item($r, 3).Value = -1482
}

# Row 64: the rank/score step up.
$ws.Cells.Item(64, 1).Value = 44368
$ws.Cells.Item(64, 2).Value = 129
$ws.Cells.Item(64, 3).Value = -1326

# Rows 65-101: dates only, no score/rank recorded yet.
for ($r = 65; $r -le 101; $r++) {
    $date = 44369 + ($r - 65)
    $ws.Cells.Item($r, 1).Value = $date
}

# Give the new date cells (A27:A101) the same date number format as the
# rest of column A, by copying the existing cell format from A26 rather
# than re-declaring a new (duplicate) number format.
$ws.Range("A26").Copy() | Out-Null
$ws.Range("A27:A101").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Update the chart's series so it keeps referencing the whole table
# (Sheet1!$A$2:$A$76 / Sheet1!$B$2:$B$76) instead of the old $26 extent.
$chartObj = $ws.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES(Sheet1!`$B`$1,Sheet1!`$A`$2:`$A`$76,Sheet1!`$B`$2:`$B`$76,1)"

# Mirror the author's final selection / cursor position.
$ws.Range("H66").Select() | Out-Null
